# Daily attendance processing - 2025-11-14 08:53:34
# Swap the "Recorded By" column (G) text from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" for the affected attendance rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,10,12,13,14,15,18,19,20,21,22,24,26,29,32,36,38,39,40,41,44,45,46,47,48,50,52,55,58,62,64,65,66,67,70,71,72,73,74,76,78,83,84,85,86,90,92,99,101,109,110,111,112,116,118,125,127,135,136,137,138,142,144,151,153)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}
